# Insert a new weekly price record for "Camote" (Peru, 2a nueva(o)) on the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Zapallo" sheet.
#
# The new record is inserted as row 165, pushing the existing rows 165-197
# down to 166-198 (dimension grows from A1:R197 to A1:R198).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 165:197 down by one row, leaving row 165 free for the new record.
$ws.Rows.Item(165).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 165 with the new weekly record.
$ws.Cells.Item(165, 1).Value  = 4
$ws.Cells.Item(165, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value  = "Los Lagos"
$ws.Cells.Item(165, 4).Value  = 44504
$ws.Cells.Item(165, 5).Value  = 10
$ws.Cells.Item(165, 6).Value  = 100112045
$ws.Cells.Item(165, 7).Value  = "Zapallo"
$ws.Cells.Item(165, 8).Value  = "Camote"
$ws.Cells.Item(165, 9).Value  = "2a nueva(o)"
$ws.Cells.Item(165, 10).Value = 150
$ws.Cells.Item(165, 11).Value = 600
$ws.Cells.Item(165, 12).Value = 600
$ws.Cells.Item(165, 13).Value = 600
$ws.Cells.Item(165, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(165, 15).Value = "Perú"
$ws.Cells.Item(165, 16).Value = 600
$ws.Cells.Item(165, 17).Value = 1
$ws.Cells.Item(165, 18).Value = "Hortaliza"
